# Update the "Periodo Mora" (period) list so it runs in ascending
# order (2106 -> 2112) instead of descending (2112 -> 2106), and move
# the 21333 "Valor Mora" amount so it stays attached to period 2112
# (now the last row) while every other row keeps 40000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2106", "2107", "2108", "2109", "2110", "2111", "2112")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

$ws.Cells.Item(16, 6).Value = 40000
$ws.Cells.Item(22, 6).Value = 21333
